$d = $word.ActiveDocument

# "Play Defenders of the Realm for Free - Review" -> "Play Defenders of the Realm for Free"
# (appears twice: the H1 title and the bold "meta title" paragraph near the end -
#  Replace:=2 / wdReplaceAll replaces every occurrence found in $d.Content in one call)
$d.Content.Find.Execute("Play Defenders of the Realm for Free - Review", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Play Defenders of the Realm for Free", 2)

# "What we like" bullet: Engaging Beat Boxes function... -> Simple and engaging gameplay
$d.Content.Find.Execute("Engaging Beat Boxes function that increases payouts", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Simple and engaging gameplay", 2)

# "What we like" bullet: Multilevel Pick Bonus feature -> Pick Bonus feature adds excitement and rewards
$d.Content.Find.Execute("Multilevel Pick Bonus feature", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pick Bonus feature adds excitement and rewards", 2)

# "What we don't like" bullet: Only 20 fixed paylines -> Limited number of paylines
$d.Content.Find.Execute("Only 20 fixed paylines", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Limited number of paylines", 2)

# "What we don't like" bullet: Bonus rounds may not trigger frequently -> No progressive jackpot feature
$d.Content.Find.Execute("Bonus rounds may not trigger frequently", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "No progressive jackpot feature", 2)

# Italic "meta description" paragraph
$d.Content.Find.Execute("Read our review of Defenders of the Realm and play for free. Engaging Beat Boxes, stunning graphics, and high volatility with a multilevel Pick Bonus.", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Read our review of Defenders of the Realm and play for free to experience the epic gameplay and stunning graphics.", 2)
